$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy R1's format (font/border) onto the two new header
#     cells T1/U1, then give all three (R1/T1/U1) the same yellow fill,
#     matching the merged style the diff shows (fontId=3, fillId=yellow,
#     borderId=2). S1 keeps its original style/text untouched.
$ws.Range("R1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("U1").PasteSpecial(-4122)

$ws.Range("T1").Value = "latitud"
$ws.Range("U1").Value = "longitud"

$ws.Range("R1").Interior.Color = 65535
$ws.Range("T1").Interior.Color = 65535
$ws.Range("U1").Interior.Color = 65535

# --- tipo_recibo column (R): rows 14-20 change from "O" to "R"
$rRows = @(14, 15, 16, 17, 18, 19, 20)
foreach ($r in $rRows) {
    $ws.Range("R$r").Value = "R"
}

# --- New latitud/longitud data columns (T, U) for every data row.
#     Copy formatting from the matching R-column data cell (style 7:
#     default font, thin border, no fill) onto T/U, then set the values.
$coords = @{
    2  = @(-11.99226322, -77.016212699999997)
    3  = @(-11.998517140000001, -77.015098710000004)
    4  = @(-11.99195128, -77.016578960000004)
    5  = @(-11.99195128, -77.016578960000004)
    6  = @(-11.99370599, -77.012533450000006)
    7  = @(-11.996716790000001, -77.015934130000005)
    8  = @(-11.996716790000001, -77.015934130000005)
    9  = @(-11.99219227, -77.0166234)
    10 = @(-11.99219227, -77.0166234)
    11 = @(-11.99219227, -77.0166234)
    12 = @(-11.99219227, -77.0166234)
    13 = @(-11.996657969999999, -77.015911099999997)
    14 = @(-11.996657969999999, -77.015911099999997)
    15 = @(-11.99194726, -77.016515279999993)
    16 = @(-11.99194726, -77.016515279999993)
    17 = @(-11.99194726, -77.016515279999993)
    18 = @(-11.991934860000001, -77.016324979999993)
    19 = @(-11.99677621, -77.015957510000007)
    20 = @(-11.99677621, -77.015957510000007)
}

for ($r = 2; $r -le 20; $r++) {
    $ws.Range("R$r").Copy()
    $ws.Range("T$r").PasteSpecial(-4122)
    $ws.Range("U$r").PasteSpecial(-4122)

    $pair = $coords[$r]
    $ws.Range("T$r").Value = $pair[0]
    $ws.Range("U$r").Value = $pair[1]
}

# --- View/selection bookkeeping to match the edited file.
$ws.Range("T1:U1").Select()
